$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '24.346.00'
$ws.Range('E2').Value = '  +9.29%  '
$ws.Range('D3').Value = '1.677.98'
$ws.Range('E3').Value = '  +4.95%  '
$ws.Range('E4').Value = '  -0.60%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '306.85'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +6.24%  '
$ws.Range('E6').Value = '  +0.05%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3711'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3438'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +1.86%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '48.19'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +13.77%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.182'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +3.66%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07251'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +2.87%  '
$ws.Range('E12').Value = '  -0.39%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '20.37'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +2.91%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.104'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +3.03%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.747'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +1.32%  '
$ws.Range('D16').Value = '1.678.41'
$ws.Range('E16').Value = '  +4.96%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001109'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +2.46%  '
$ws.Range('E18').Value = '  +0.16%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06724'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +1.31%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '81.13'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +3.78%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '16.44'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +1.59%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.089'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +0.83%  '
$ws.Range('E23').Value = '  +1.49%  '
$ws.Range('D24').Value = '24.314.67'
$ws.Range('E24').Value = '  +8.87%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.427'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +1.28%  '
$ws.Range('E26').Value = '  -12.08%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.659'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +6.68%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '152.23'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +0.80%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '19.56'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -0.49%  '
$ws.Range('D30').Value = '1.863.05'
$ws.Range('E30').Value = '  +4.76%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '127.21'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +5.36%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.309'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +5.47%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.036'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -3.60%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.9671'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +2.23%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.738'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +7.81%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.08473'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +2.59%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '8.995'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +4.04%  '
$ws.Range('E38').Value = '  +4.12%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.06421'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +4.47%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.341'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +0.83%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.02331'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +5.47%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.261'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +1.87%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.2107'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +4.18%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.6163'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +3.92%  '
$ws.Range('E45').Value = '  +0.24%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.778'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +3.29%  '
$ws.Range('B47').Value = 'EnergySwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '12.98'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -1.78%  '
$ws.Range('B48').Value = 'Decentraland'
$ws.Range('C48').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.5937'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +3.80%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '126.96'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +1.60%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.025'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +2.76%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.07207'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +5.60%  '
